# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" (Overview sheet) and the
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" (per-locale
# sheets) for the file "0428f7e4-a6e9-420e-8672-b3d150ae611e.md" to reflect a
# fresh handback report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 2 = 0428f7e4-a6e9-420e-8672-b3d150ae611e.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-11-14 06:25:45"

# --- zh-cn sheet: row 2 = 0428f7e4-a6e9-420e-8672-b3d150ae611e.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-11-14 06:25:31"
$zhcn.Range("K2").Value = "2016-11-14 06:26:25"

# --- de-de sheet: row 2 = 0428f7e4-a6e9-420e-8672-b3d150ae611e.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-11-14 06:25:45"
$dede.Range("K2").Value = "2016-11-14 06:26:44"
